# A new daily price record for "Pepino ensalada" at Terminal La Palmera de
# La Serena needs to be inserted as row 847 (sheet is sorted by category and
# the new observation is dated 2023-08-04 / Excel serial 45142). Inserting a
# whole row there pushes the existing rows 847:943 down to 848:944, which is
# exactly what the target workbook shows (dimension grows from R943 to R944
# and every subsequent row's data shifts down by one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 847; everything from 847 downward shifts to 848+.
$ws.Rows.Item(847).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(847, 1).Value  = 8
$ws.Cells.Item(847, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(847, 3).Value  = "Coquimbo"
$ws.Cells.Item(847, 4).Value  = 45142
$ws.Cells.Item(847, 5).Value  = 4
$ws.Cells.Item(847, 6).Value  = 100112043
$ws.Cells.Item(847, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(847, 8).Value  = "Sin especificar"
$ws.Cells.Item(847, 9).Value  = "Primera"
$ws.Cells.Item(847, 10).Value = 400
$ws.Cells.Item(847, 11).Value = 10000
$ws.Cells.Item(847, 12).Value = 11000
$ws.Cells.Item(847, 13).Value = 10500
$ws.Cells.Item(847, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(847, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(847, 16).Value = 175
$ws.Cells.Item(847, 17).Value = 60
$ws.Cells.Item(847, 18).Value = "Hortaliza"
